$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.509.88"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.571.27"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.25"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0591"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0883"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "1.795.03"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "1.579.05"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").Value = "28.486.45"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "226.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("E23").Value = "  -5.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.54%  "
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "1.398.00"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("E39").Value = "  +3.42%  "
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.535"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.979"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("D48").Value = "1.707.49"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").Value = "  -4.06%  "
$ws.Range("E51").Value = "  -1.44%  "
